# 977: Add GS to extract process and GS tab to example files
$wb = $excel.ActiveWorkbook

# Add the new "GS" worksheet after the last existing sheet (CMS), so it
# becomes the new final/active tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "GS"

# Header row for the new GS sheet.
$headers = @("Contact_ID", "Contact_Date", "Contact_Type_Code", "Contact_Type_Desc", "OM_Name", "OM_Key", "OM_Grade", "OM_Team_Key", "OM_Provider_Code")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Color = 0
}

$ws.Range("A1:I1").Select()
